# Apply crypto price/volume updates from the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.730.21'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').Value = '3.850.12'
$ws.Range('E3').Value = '  -1.64%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.59'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.20'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('D7').Value = '3.846.80'
$ws.Range('E7').Value = '  -1.71%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.526'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.42%  '
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.34'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.455'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.83'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.35%  '
$ws.Range('D15').Value = '4.498.04'
$ws.Range('E15').Value = '  -1.55%  '
$ws.Range('D16').Value = '3.858.90'
$ws.Range('E16').Value = '  -1.34%  '
$ws.Range('D17').Value = '67.826.47'
$ws.Range('E17').Value = '  -1.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.04'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +5.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.33'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('E20').Value = '  -1.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.92'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '462.79'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -4.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.729'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.56%  '
$ws.Range('E24').Value = '  -3.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.11'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.25'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.12'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.92%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.99'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.76%  '
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('D31').Value = '3.999.06'
$ws.Range('E31').Value = '  -1.65%  '
$ws.Range('E32').Value = '  -1.15%  '
$ws.Range('E33').Value = '  -2.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '31.04'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.29'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.67%  '
$ws.Range('D36').Value = '3.825.63'
$ws.Range('E36').Value = '  -0.90%  '
$ws.Range('E37').Value = '  -2.10%  '
$ws.Range('E38').Value = '  -1.64%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.23'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +5.07%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.310'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '427.06'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.96'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.61%  '
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '47.16'
$ws.Range('D47').ClearFormats()
$ws.Range('E48').Value = '  +0.70%  '
$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '40.65'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.53%  '
$ws.Range('B50').Value = 'FLOKI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000273'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '143.86'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.41%  '
